# Update data from Streamlit app
# Row 9 corresponds to the "Oppo Service Center" lease entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lease End Date: 45822 -> 45091
$ws.Range("C9").Value = 45091

# Lease Duration (Years): 3 -> 1
$ws.Range("E9").Value = 1

# Actual Income (Rp/year): 75000000 -> 81000000
$ws.Range("G9").Value = 81000000

# Payment Scheme: "Split Per Year" -> "Full Lease Upfront"
$ws.Range("H9").Value = "Full Lease Upfront"
